$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Log a new work entry in row 45: date, hours, and description
$ws.Range("A45").Value = 44074
$ws.Range("B45").Value = 7
$ws.Range("C45").Value = "Html sivujen ja css:n toteuttamista"

# Match the saved scroll position / selection of the sheet view
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C45").Select()
